$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = -0.02080121282546088
$ws.Cells.Item(1, 2).Value = 0.020769208692158259
$ws.Cells.Item(2, 1).Value = 0.0274867939465544
$ws.Cells.Item(2, 2).Value = -0.027702376910264448
$ws.Cells.Item(3, 1).Value = 0.13062804871822919
$ws.Cells.Item(3, 2).Value = -0.13103475602143178
$ws.Cells.Item(4, 1).Value = -0.15295002207086839
$ws.Cells.Item(4, 2).Value = 0.15238838180301428
$ws.Cells.Item(5, 1).Value = -0.1463883826931589
$ws.Cells.Item(5, 2).Value = 0.14527243770483
$ws.Cells.Item(6, 1).Value = -0.049229453972952086
$ws.Cells.Item(6, 2).Value = 0.049200610167135395
$ws.Cells.Item(7, 1).Value = -0.02920061122884654
$ws.Cells.Item(7, 2).Value = 0.029161621827936912
$ws.Cells.Item(8, 1).Value = -0.0091616228934734423
$ws.Cells.Item(8, 2).Value = 0.0091378156264765309
$ws.Cells.Item(9, 1).Value = -0.0031378165494979626
$ws.Cells.Item(9, 2).Value = 0.0031138033113640162
$ws.Cells.Item(10, 1).Value = 0.0028861957648089742
$ws.Cells.Item(10, 2).Value = -0.0028855518268500191
$ws.Cells.Item(11, 1).Value = 0.0073855509184674872
$ws.Cells.Item(11, 2).Value = -0.007396728263810104
$ws.Cells.Item(12, 1).Value = -0.074578058519509938
$ws.Cells.Item(12, 2).Value = 0.074055802554104488
$ws.Cells.Item(13, 1).Value = -0.068055803501035683
$ws.Cells.Item(13, 2).Value = 0.067892031880574244
$ws.Cells.Item(14, 1).Value = -0.055892032897072674
$ws.Cells.Item(14, 2).Value = 0.055763153095154117
$ws.Cells.Item(15, 1).Value = -0.049763154055975534
$ws.Cells.Item(15, 2).Value = 0.04962905994780975
$ws.Cells.Item(16, 1).Value = -0.015026058353548422
$ws.Cells.Item(16, 2).Value = 0.015003807062009322
$ws.Cells.Item(17, 1).Value = -0.0090038080335990145
$ws.Cells.Item(17, 2).Value = 0.0089999989967237681
$ws.Cells.Item(18, 1).Value = -0.038256003569525632
$ws.Cells.Item(18, 2).Value = 0.038245279217154149
$ws.Cells.Item(19, 1).Value = -0.029245280138309848
$ws.Cells.Item(19, 2).Value = 0.029164761032025144
$ws.Cells.Item(20, 1).Value = -0.01801284919202395
$ws.Cells.Item(20, 2).Value = 0.018004144316890347
$ws.Cells.Item(21, 1).Value = -0.0090041452474975969
$ws.Cells.Item(21, 2).Value = 0.0089999990686502329
$ws.Cells.Item(22, 1).Value = -0.093925672726953735
$ws.Cells.Item(22, 2).Value = 0.093620451275262084
$ws.Cells.Item(23, 1).Value = -0.084620452201149554
$ws.Cells.Item(23, 2).Value = 0.084123770382761975
$ws.Cells.Item(24, 1).Value = -0.042123771664473786
$ws.Cells.Item(24, 2).Value = 0.041999998711357733
$ws.Cells.Item(25, 1).Value = -0.090042984647382696
$ws.Cells.Item(25, 2).Value = 0.089943948341133506
$ws.Cells.Item(26, 1).Value = -0.083943949259747797
$ws.Cells.Item(26, 2).Value = 0.083819488147916843
$ws.Cells.Item(27, 1).Value = -0.077819489070630077
$ws.Cells.Item(27, 2).Value = 0.077405678593447824
$ws.Cells.Item(28, 1).Value = -0.071405679534485067
$ws.Cells.Item(28, 2).Value = 0.071139601951500708
$ws.Cells.Item(29, 1).Value = -0.059139602966679305
$ws.Cells.Item(29, 2).Value = 0.059026081616801562
$ws.Cells.Item(30, 1).Value = -0.039026082720456046
$ws.Cells.Item(30, 2).Value = 0.03867829590085714
$ws.Cells.Item(31, 1).Value = -0.035690404934687336
$ws.Cells.Item(31, 2).Value = 0.035665321323333998
$ws.Cells.Item(32, 1).Value = -0.0060005686515998136
$ws.Cells.Item(32, 2).Value = 0.00599999903044246
